{"js": "// Replace the 25 division fact-family answer strings in the document's\n// single table with their new values, matching the exact text (including\n// the division sign \"\u00f7\" and trailing \", remainder\" suffix) so each\n// replacement only touches its own cell.\nconst replacements = [\n  [\"25\u00f75=5, 0\", \"69\u00f79=7, 6\"],\n  [\"72\u00f76=12, 0\", \"91\u00f72=45, 1\"],\n  [\"38\u00f78=4, 6\", \"61\u00f77=8, 5\"],\n  [\"98\u00f76=16, 2\", \"68\u00f79=7, 5\"],\n  [\"61\u00f74=15, 1\", \"53\u00f75=10, 3\"],\n  [\"74\u00f74=18, 2\", \"71\u00f77=10, 1\"],\n  [\"73\u00f79=8, 1\", \"32\u00f75=6, 2\"],\n  [\"28\u00f72=14, 0\", \"64\u00f77=9, 1\"],\n  [\"62\u00f76=10, 2\", \"73\u00f76=12, 1\"],\n  [\"35\u00f76=5, 5\", \"20\u00f77=2, 6\"],\n  [\"60\u00f76=10, 0\", \"27\u00f76=4, 3\"],\n  [\"79\u00f74=19, 3\", \"36\u00f76=6, 0\"],\n  [\"54\u00f72=27, 0\", \"42\u00f77=6, 0\"],\n  [\"71\u00f74=17, 3\", \"66\u00f79=7, 3\"],\n  [\"80\u00f79=8, 8\", \"51\u00f77=7, 2\"],\n  [\"81\u00f79=9, 0\", \"14\u00f74=3, 2\"],\n  [\"85\u00f72=42, 1\", \"26\u00f74=6, 2\"],\n  [\"50\u00f79=5, 5\", \"44\u00f74=11, 0\"],\n  [\"24\u00f73=8, 0\", \"76\u00f73=25, 1\"],\n  [\"32\u00f76=5, 2\", \"26\u00f79=2, 8\"],\n  [\"85\u00f79=9, 4\", \"43\u00f79=4, 7\"],\n  [\"61\u00f73=20, 1\", \"46\u00f77=6, 4\"],\n  [\"20\u00f75=4, 0\", \"62\u00f72=31, 0\"],\n  [\"84\u00f76=14, 0\", \"38\u00f74=9, 2\"],\n  [\"83\u00f75=16, 3\", \"48\u00f74=12, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  // matchCase keeps each lookup anchored to its exact, unique original\n  // string, so every cell is updated independently in a single pass.\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Old answer -> new answer for each of the 25 filled table cells.\n$pairs = @(\n    @(\"25\u00f75=5, 0\", \"69\u00f79=7, 6\"),\n    @(\"72\u00f76=12, 0\", \"91\u00f72=45, 1\"),\n    @(\"38\u00f78=4, 6\", \"61\u00f77=8, 5\"),\n    @(\"98\u00f76=16, 2\", \"68\u00f79=7, 5\"),\n    @(\"61\u00f74=15, 1\", \"53\u00f75=10, 3\"),\n    @(\"74\u00f74=18, 2\", \"71\u00f77=10, 1\"),\n    @(\"73\u00f79=8, 1\", \"32\u00f75=6, 2\"),\n    @(\"28\u00f72=14, 0\", \"64\u00f77=9, 1\"),\n    @(\"62\u00f76=10, 2\", \"73\u00f76=12, 1\"),\n    @(\"35\u00f76=5, 5\", \"20\u00f77=2, 6\"),\n    @(\"60\u00f76=10, 0\", \"27\u00f76=4, 3\"),\n    @(\"79\u00f74=19, 3\", \"36\u00f76=6, 0\"),\n    @(\"54\u00f72=27, 0\", \"42\u00f77=6, 0\"),\n    @(\"71\u00f74=17, 3\", \"66\u00f79=7, 3\"),\n    @(\"80\u00f79=8, 8\", \"51\u00f77=7, 2\"),\n    @(\"81\u00f79=9, 0\", \"14\u00f74=3, 2\"),\n    @(\"85\u00f72=42, 1\", \"26\u00f74=6, 2\"),\n    @(\"50\u00f79=5, 5\", \"44\u00f74=11, 0\"),\n    @(\"24\u00f73=8, 0\", \"76\u00f73=25, 1\"),\n    @(\"32\u00f76=5, 2\", \"26\u00f79=2, 8\"),\n    @(\"85\u00f79=9, 4\", \"43\u00f79=4, 7\"),\n    @(\"61\u00f73=20, 1\", \"46\u00f77=6, 4\"),\n    @(\"20\u00f75=4, 0\", \"62\u00f72=31, 0\"),\n    @(\"84\u00f76=14, 0\", \"38\u00f74=9, 2\"),\n    @(\"83\u00f75=16, 3\", \"48\u00f74=12, 0\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    # Each old value is a unique literal string, so a fresh Find scoped to\n    # the whole document body safely retargets just its own cell.\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n    #         MatchSoundsLike, MatchAllWordForms, Forward, Wrap,\n    #         Format, ReplaceWith, Replace=wdReplaceAll)\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
